$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 454, shifting existing rows 454:509 down to 455:510
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new record's data
$ws.Range("A454").Value = 4
$ws.Range("B454").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C454").Value = "Los Lagos"
$ws.Range("D454").Value = 45142
$ws.Range("D454").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E454").Value = 10
$ws.Range("F454").Value = 100112017
$ws.Range("G454").Value = "Apio"
$ws.Range("H454").Value = "Americana (o)"
$ws.Range("I454").Value = "Primera"
$ws.Range("J454").Value = 45
$ws.Range("K454").Value = 11000
$ws.Range("L454").Value = 11000
$ws.Range("M454").Value = 11000
$ws.Range("N454").Value = "`$/docena de matas"
$ws.Range("O454").Value = "Región de Coquimbo"
$ws.Range("P454").Value = 1833
$ws.Range("Q454").Value = 6
$ws.Range("R454").Value = "Hortaliza"
